$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5709.4443
$ws.Range("I51").Value = 7582.75
$ws.Range("J51").Value = 4210.8
$ws.Range("K51").Value = 7582.75
$ws.Range("L51").Value = 4210.8
$ws.Range("M51").Value = -7098.75
$ws.Range("N51").Value = -5178.8
$ws.Range("H96").Value = 4867.9
$ws.Range("I96").Value = 5522.375
$ws.Range("K96").Value = 16567.125
$ws.Range("M96").Value = -15194.125
$ws.Range("H99").Value = 836.7646999999999
$ws.Range("I99").Value = 2588
$ws.Range("J99").Value = 297.92307
$ws.Range("K99").Value = 7764
$ws.Range("L99").Value = 893.7692099999999
$ws.Range("M99").Value = -6266
$ws.Range("N99").Value = -3889.76921
$ws.Range("H106").Value = 22224948
$ws.Range("I106").Value = 30305612
$ws.Range("J106").Value = 3125
$ws.Range("K106").Value = 30305612
$ws.Range("L106").Value = 3125
$ws.Range("M106").Value = -30304981
$ws.Range("N106").Value = -4387
$ws.Range("H132").Value = 4920.7256
$ws.Range("I132").Value = 2543
$ws.Range("K132").Value = 7629
$ws.Range("M132").Value = -5099
$ws.Range("H137").Value = 7182.5
$ws.Range("I137").Value = 8578.385
$ws.Range("J137").Value = 1133.6666
$ws.Range("K137").Value = 25735.155
$ws.Range("L137").Value = 3400.9998
$ws.Range("M137").Value = -23185.155
$ws.Range("N137").Value = -8500.9998
$ws.Range("H138").Value = 2335.3699
$ws.Range("I138").Value = 1881.3928
$ws.Range("K138").Value = 5644.178400000001
$ws.Range("M138").Value = -504.1784000000007

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 350.66666
$ws.Range("I5").Value = 334
$ws.Range("K5").Value = 334
$ws.Range("M5").Value = -222
$ws.Range("H61").Value = 5544.7617
$ws.Range("I61").Value = 8862.143
$ws.Range("J61").Value = 3886.0715
$ws.Range("K61").Value = 8862.143
$ws.Range("L61").Value = 3886.0715
$ws.Range("M61").Value = -8650.143
$ws.Range("N61").Value = -4310.0715
$ws.Range("H110").Value = 1749.9131
$ws.Range("I110").Value = 1829.6097
$ws.Range("J110").Value = 1096.4
$ws.Range("K110").Value = 1829.6097
$ws.Range("L110").Value = 1096.4
$ws.Range("M110").Value = 215.3903
$ws.Range("N110").Value = -5186.4
$ws.Range("H122").Value = 1823.3658
$ws.Range("I122").Value = 1932.7778
$ws.Range("J122").Value = 1612.3572
$ws.Range("K122").Value = 5798.3334
$ws.Range("L122").Value = 4837.071599999999
$ws.Range("M122").Value = -3348.3334
$ws.Range("N122").Value = -9737.071599999999
$ws.Range("H132").Value = 2769.2
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 5544.7617
$ws.Range("I136").Value = 8862.143
$ws.Range("J136").Value = 3886.0715
$ws.Range("K136").Value = 26586.429
$ws.Range("L136").Value = 11658.2145
$ws.Range("M136").Value = -24036.429
$ws.Range("N136").Value = -16758.2145
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 350.66666
$ws.Range("I4").Value = 334
$ws.Range("K4").Value = 334
$ws.Range("M4").Value = -219
$ws.Range("H5").Value = 1538.7778
$ws.Range("I5").Value = 1516.3334
$ws.Range("J5").Value = 1550
$ws.Range("K5").Value = 1516.3334
$ws.Range("L5").Value = 1550
$ws.Range("M5").Value = -1403.3334
$ws.Range("N5").Value = -1776
$ws.Range("H10").Value = 3001.5
$ws.Range("J10").Value = 3668.6667
$ws.Range("L10").Value = 3668.6667
$ws.Range("N10").Value = -3948.6667
$ws.Range("H11").Value = 728.8570999999999
$ws.Range("I11").Value = 433.66666
$ws.Range("K11").Value = 433.66666
$ws.Range("M11").Value = -293.66666
$ws.Range("H12").Value = 1500
$ws.Range("I12").Value = 1500
$ws.Range("K12").Value = 1500
$ws.Range("M12").Value = -1332
$ws.Range("H16").Value = 11000
$ws.Range("I16").Value = 10000
$ws.Range("J16").Value = 12000
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = -9830
$ws.Range("N16").Value = -12340
$ws.Range("H17").Value = 55555
$ws.Range("I17").Value = 55555
$ws.Range("K17").Value = 55555
$ws.Range("M17").Value = -55383
$ws.Range("H19").Value = 6620
$ws.Range("I19").Value = 6620
$ws.Range("K19").Value = 6620
$ws.Range("M19").Value = -6447
$ws.Range("H105").Value = 1758.7567
$ws.Range("I105").Value = 1885.2273
$ws.Range("K105").Value = 1885.2273
$ws.Range("M105").Value = -138.2273
$ws.Range("H107").Value = 2284.2856
$ws.Range("I107").Value = 2019.375
$ws.Range("J107").Value = 3132
$ws.Range("K107").Value = 2019.375
$ws.Range("L107").Value = 3132
$ws.Range("M107").Value = -99.375
$ws.Range("N107").Value = -6972
$ws.Range("H134").Value = 2086.5
$ws.Range("J134").Value = 1198.5
$ws.Range("L134").Value = 3595.5
$ws.Range("N134").Value = -8665.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 9195.846
$ws.Range("I22").Value = 14784.875
$ws.Range("K22").Value = 14784.875
$ws.Range("M22").Value = -14434.875
$ws.Range("H31").Value = 3076.3242
$ws.Range("I31").Value = 1755.4546
$ws.Range("J31").Value = 5013.6
$ws.Range("K31").Value = 1755.4546
$ws.Range("L31").Value = 5013.6
$ws.Range("M31").Value = -1460.4546
$ws.Range("N31").Value = -5603.6
$ws.Range("H34").Value = 3076.3242
$ws.Range("I34").Value = 1755.4546
$ws.Range("J34").Value = 5013.6
$ws.Range("K34").Value = 1755.4546
$ws.Range("L34").Value = 5013.6
$ws.Range("M34").Value = -1553.4546
$ws.Range("N34").Value = -5417.6
$ws.Range("H105").Value = 1045.44
$ws.Range("I105").Value = 953.4211
$ws.Range("K105").Value = 953.4211
$ws.Range("M105").Value = 793.5789

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 50
$ws.Range("I41").Value = 50
$ws.Range("K41").Value = 150
$ws.Range("M41").Value = 188
$ws.Range("H92").Value = 759.1667
$ws.Range("I92").Value = 704.7143
$ws.Range("J92").Value = 835.4
$ws.Range("K92").Value = 2114.1429
$ws.Range("L92").Value = 2506.2
$ws.Range("M92").Value = -866.1428999999998
$ws.Range("N92").Value = -5002.2
$ws.Range("H99").Value = 2346
$ws.Range("I99").Value = 1169.6
$ws.Range("K99").Value = 3508.8
$ws.Range("M99").Value = -1262.8
$ws.Range("H129").Value = 1066.5
$ws.Range("I129").Value = 759.2857
$ws.Range("J129").Value = 1496.6
$ws.Range("K129").Value = 2277.8571
$ws.Range("L129").Value = 4489.799999999999
$ws.Range("M129").Value = 2722.1429
$ws.Range("N129").Value = -14489.8
$ws.Range("H131").Value = 2102682.8
$ws.Range("I131").Value = 3677274
$ws.Range("K131").Value = 11031822
$ws.Range("M131").Value = -11026782
$ws.Range("H132").Value = 2847.875
$ws.Range("I132").Value = 1186
$ws.Range("K132").Value = 10674
$ws.Range("M132").Value = -8144

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2936.2666
$ws.Range("I46").Value = 1078.7142
$ws.Range("J46").Value = 4561.625
$ws.Range("K46").Value = 1078.7142
$ws.Range("L46").Value = 4561.625
$ws.Range("M46").Value = -890.7141999999999
$ws.Range("N46").Value = -4937.625
$ws.Range("H61").Value = 1589.409
$ws.Range("I61").Value = 1648.85
$ws.Range("J61").Value = 995
$ws.Range("K61").Value = 1648.85
$ws.Range("L61").Value = 995
$ws.Range("M61").Value = -1446.85
$ws.Range("N61").Value = -1399
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 1589.409
$ws.Range("I113").Value = 1648.85
$ws.Range("J113").Value = 995
$ws.Range("K113").Value = 1648.85
$ws.Range("L113").Value = 995
$ws.Range("M113").Value = 521.1500000000001
$ws.Range("N113").Value = -5335
$ws.Range("H132").Value = 18372.854
$ws.Range("I132").Value = 22255.473
$ws.Range("J132").Value = 4654.2666
$ws.Range("K132").Value = 66766.41900000001
$ws.Range("L132").Value = 13962.7998
$ws.Range("M132").Value = -64236.41900000001
$ws.Range("N132").Value = -19022.7998

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1454.7778
$ws.Range("I100").Value = 1298.8572
$ws.Range("K100").Value = 2597.7144
$ws.Range("M100").Value = -2056.7144
$ws.Range("H107").Value = 507.75
$ws.Range("I107").Value = 514.3
$ws.Range("J107").Value = 475
$ws.Range("K107").Value = 1542.9
$ws.Range("L107").Value = 1425
$ws.Range("M107").Value = 377.1000000000001
$ws.Range("N107").Value = -5265
$ws.Range("H132").Value = 7846.4053
$ws.Range("I132").Value = 7796.75
$ws.Range("J132").Value = 8164.2
$ws.Range("K132").Value = 23390.25
$ws.Range("L132").Value = 24492.6
$ws.Range("M132").Value = -20860.25
$ws.Range("N132").Value = -29552.6
